$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 452.25
$ws.Range("I12").Value = 319.83334
$ws.Range("K12").Value = 319.83334
$ws.Range("M12").Value = -149.83334

$ws.Range("H33").Value = 860.5
$ws.Range("I33").Value = 827.4286
$ws.Range("K33").Value = 827.4286
$ws.Range("M33").Value = -598.4286

$ws.Range("H41").Value = 436.72726
$ws.Range("I41").Value = 440.33334
$ws.Range("K41").Value = 440.33334
$ws.Range("M41").Value = -0.3333400000000211

$ws.Range("H98").Value = 1390.6552
$ws.Range("I98").Value = 1510.6364
$ws.Range("K98").Value = 1510.6364
$ws.Range("M98").Value = -12.63640000000009

$ws.Range("H122").Value = 1390.6552
$ws.Range("I122").Value = 1510.6364
$ws.Range("K122").Value = 4531.9092
$ws.Range("M122").Value = -2081.9092

$ws.Range("H132").Value = 994037
$ws.Range("I132").Value = 1074790.1
$ws.Range("K132").Value = 3224370.3
$ws.Range("M132").Value = -3221840.3

$ws.Range("H136").Value = 148988.5
$ws.Range("J136").Value = 148988.5
$ws.Range("L136").Value = 148988.5
$ws.Range("N136").Value = -159188.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1097.3334
$ws.Range("I4").Value = 895
$ws.Range("J4").Value = 1502
$ws.Range("K4").Value = 895
$ws.Range("L4").Value = 1502
$ws.Range("M4").Value = -779
$ws.Range("N4").Value = -1734

$ws.Range("H74").Value = 5576.304
$ws.Range("I74").Value = 5311.357
$ws.Range("J74").Value = 5988.4443
$ws.Range("K74").Value = 5311.357
$ws.Range("L74").Value = 5988.4443
$ws.Range("M74").Value = -4437.357
$ws.Range("N74").Value = -7736.4443

$ws.Range("H77").Value = 5576.304
$ws.Range("I77").Value = 5311.357
$ws.Range("J77").Value = 5988.4443
$ws.Range("K77").Value = 26556.785
$ws.Range("L77").Value = 29942.2215
$ws.Range("M77").Value = -22188.785
$ws.Range("N77").Value = -38678.2215

$ws.Range("H110").Value = 4517.788
$ws.Range("I110").Value = 3255.0588
$ws.Range("K110").Value = 3255.0588
$ws.Range("M110").Value = -1210.0588

$ws.Range("H133").Value = 179666
$ws.Range("J133").Value = 179666
$ws.Range("L133").Value = 179666
$ws.Range("N133").Value = -184726

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 490
$ws.Range("J10").Value = 490
$ws.Range("L10").Value = 490
$ws.Range("N10").Value = -770

$ws.Range("H11").Value = 2734.1667
$ws.Range("I11").Value = 2880
$ws.Range("J11").Value = 2005
$ws.Range("K11").Value = 2880
$ws.Range("L11").Value = 2005
$ws.Range("M11").Value = -2740
$ws.Range("N11").Value = -2285

$ws.Range("H81").Value = 16120
$ws.Range("J81").Value = 16120
$ws.Range("L81").Value = 16120
$ws.Range("N81").Value = -18242

$ws.Range("H84").Value = 16120
$ws.Range("J84").Value = 16120
$ws.Range("L84").Value = 48360
$ws.Range("N84").Value = -58968

$ws.Range("H99").Value = 5033.12
$ws.Range("I99").Value = 2460.9412
$ws.Range("K99").Value = 2460.9412
$ws.Range("M99").Value = -962.9412000000002

$ws.Range("H134").Value = 10611.95
$ws.Range("I134").Value = 10217
$ws.Range("J134").Value = 12850
$ws.Range("K134").Value = 30651
$ws.Range("L134").Value = 38550
$ws.Range("M134").Value = -28116
$ws.Range("N134").Value = -43620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8562955
$ws.Range("I4").Value = 11000442
$ws.Range("J4").Value = 437996.34
$ws.Range("K4").Value = 33001326
$ws.Range("L4").Value = 1313989.02
$ws.Range("M4").Value = -33001214
$ws.Range("N4").Value = -1314213.02

$ws.Range("H80").Value = 2347.5
$ws.Range("J80").Value = 2347.5
$ws.Range("L80").Value = 7042.5
$ws.Range("N80").Value = -8914.5

$ws.Range("H83").Value = 2347.5
$ws.Range("J83").Value = 2347.5
$ws.Range("L83").Value = 21127.5
$ws.Range("N83").Value = -30487.5

$ws.Range("H86").Value = 287.4
$ws.Range("J86").Value = 287.4
$ws.Range("L86").Value = 862.1999999999999
$ws.Range("N86").Value = -3234.2

$ws.Range("H89").Value = 287.4
$ws.Range("J89").Value = 287.4
$ws.Range("L89").Value = 2586.6
$ws.Range("N89").Value = -14442.6

$ws.Range("H98").Value = 1771.2858
$ws.Range("I98").Value = 2333.6667
$ws.Range("J98").Value = 1349.5
$ws.Range("K98").Value = 7001.000100000001
$ws.Range("L98").Value = 4048.5
$ws.Range("M98").Value = -5503.000100000001
$ws.Range("N98").Value = -7044.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14764.143
$ws.Range("I61").Value = 12266.5
$ws.Range("K61").Value = 12266.5
$ws.Range("M61").Value = -12064.5

$ws.Range("H93").Value = 1416.579
$ws.Range("I93").Value = 2354.5
$ws.Range("K93").Value = 2354.5
$ws.Range("M93").Value = -1106.5

$ws.Range("H100").Value = 1886.7
$ws.Range("I100").Value = 1623.75
$ws.Range("K100").Value = 1623.75
$ws.Range("M100").Value = -1082.75

$ws.Range("H113").Value = 14764.143
$ws.Range("I113").Value = 12266.5
$ws.Range("K113").Value = 12266.5
$ws.Range("M113").Value = -10096.5

$ws.Range("H136").Value = 4821.306
$ws.Range("I136").Value = 5916.643
$ws.Range("J136").Value = 4383.1714
$ws.Range("K136").Value = 17749.929
$ws.Range("L136").Value = 13149.5142
$ws.Range("M136").Value = -15199.929
$ws.Range("N136").Value = -18249.5142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H62").Value = 17163
$ws.Range("I62").Value = 11648.571
$ws.Range("J62").Value = 22677.428
$ws.Range("K62").Value = 11648.571
$ws.Range("L62").Value = 22677.428
$ws.Range("M62").Value = -11024.571
$ws.Range("N62").Value = -23925.428

$ws.Range("H65").Value = 17163
$ws.Range("I65").Value = 11648.571
$ws.Range("J65").Value = 22677.428
$ws.Range("K65").Value = 58242.855
$ws.Range("L65").Value = 113387.14
$ws.Range("M65").Value = -55122.855
$ws.Range("N65").Value = -119627.14

$ws.Range("H81").Value = 1842.9231
$ws.Range("I81").Value = 1359.8182
$ws.Range("J81").Value = 4500
$ws.Range("K81").Value = 2719.6364
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -1658.6364
$ws.Range("N81").Value = -11122

$ws.Range("H84").Value = 1842.9231
$ws.Range("I84").Value = 1359.8182
$ws.Range("J84").Value = 4500
$ws.Range("K84").Value = 13598.182
$ws.Range("L84").Value = 45000
$ws.Range("M84").Value = -8294.181999999999
$ws.Range("N84").Value = -55608

$ws.Range("H96").Value = 3235
$ws.Range("I96").Value = 2083
$ws.Range("K96").Value = 2083
$ws.Range("M96").Value = -710
